# "Arreglos sugerencias de Ruben"
#
# The age-bracket breakdown on sheet "1_18" is consolidated from four
# brackets (De 0 a 4 / De 5 a 9 / De 10 a 14 / De 15 a 17) plus the
# existing "De 0 a 17" / "Entre 18 y 19" rows down to three rows:
#   De 0 a 17, De 0 a 9, De 10 a 17 - with refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1_18")

# Make this the active sheet/tab (matches the workbook's new activeTab).
$ws.Activate() | Out-Null

# Overwrite the first three data rows with the consolidated categories
# and their updated values.
$ws.Range("A2").Value = "De 0 a 17"
$ws.Range("B2").Value = 68.2
$ws.Range("A3").Value = "De 0 a 9"
$ws.Range("B3").Value = 70.2
$ws.Range("A4").Value = "De 10 a 17"
$ws.Range("B4").Value = 65.9

# Remove the now-obsolete rows that held "De 15 a 17", "De 0 a 17" and
# "Entre 18 y 19" (the last two categories are superseded by the new
# consolidated rows above), shifting everything below up.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

# Select the new header/data range, matching the updated selection.
$ws.Range("A2:B2").Select() | Out-Null
